# refactor: Enhance reusability of code and completed 'VerifyEntriesInTable' logic
#
# Inserts a new "Id" column at the front of the UserDetails sheet, numbering
# the two existing rows (1, 2), and refreshes the (former) UserName column
# values. Also updates the current selection to match the author's last
# position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A - everything currently in A:H shifts to B:I
$ws.Columns("A:A").Insert()

# The engine doesn't auto-shift the worksheet's Hyperlinks collection when a
# column is inserted, so re-anchor them to their new location (G -> H) by
# hand. Re-adding resets the cell text/style as a side effect, so restore
# both afterwards.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:admin@mail.com", "", "", "mailto:admin@mail.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:customer@mail.com", "", "", "mailto:customer@mail.com")
$ws.Range("H2").Value = "admin@mail.com"
$ws.Range("H3").Value = "customer@mail.com"
$ws.Range("H2:H3").Style = "Hyperlink"

# New "Id" header + row numbers (stored as real numbers, but displayed via a
# text ("@") number format - matches the rest of the sheet's text columns)
$ws.Range("A1").Value = "Id"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A2:A3").NumberFormat = "@"

# Refresh the UserName values (now column D after the insert). The column's
# Text ("@") format would otherwise coerce a plain ".Value = <number>" write
# into a shared string, so reset to "Normal" first to land a real number,
# then copy the text format back from the untouched header cell (D1) rather
# than toggling NumberFormat directly, which would otherwise register a
# stray custom/duplicate style.
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").Value = 1351513804
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").Value = 8318442633

$ws.Range("D1").Copy()
$ws.Range("D2:D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the author's final cursor position
$ws.Range("E9").Select()
